# Scheduled runner update: refresh computed market-price / leve-profit figures
# across the per-job "Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Cells hold plain (non-formula) numeric snapshots pulled from the market-board
# API on each run, so this script simply rewrites the affected cells with
# their latest values. A handful of profit cells (column M/N) flip between a
# blank cell and a real number from run to run (e.g. when a price swings from
# a loss to break-even or vice versa) - those are written as $null to clear
# the cell entirely, matching how the source data omits them.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 117.333336
$ws.Range("I9").Value = 42.833332
$ws.Range("K9").Value = 42.833332
$ws.Range("M9").Value = 126.166668
$ws.Range("H11").Value = 577.6667
$ws.Range("I11").Value = 577.6667
$ws.Range("K11").Value = 577.6667
$ws.Range("M11").Value = -437.6667
$ws.Range("H17").Value = 4500.1665
$ws.Range("J17").Value = 6000
$ws.Range("L17").Value = 18000
$ws.Range("N17").Value = -18336
$ws.Range("H62").Value = 4943.9
$ws.Range("I62").Value = 4349.6
$ws.Range("K62").Value = 4349.6
$ws.Range("M62").Value = -3725.6
$ws.Range("H65").Value = 4943.9
$ws.Range("I65").Value = 4349.6
$ws.Range("K65").Value = 21748
$ws.Range("M65").Value = -18628
$ws.Range("H106").Value = 3633
$ws.Range("I106").Value = 3633
$ws.Range("K106").Value = 3633
$ws.Range("M106").Value = -3002
$ws.Range("H107").Value = 1934.0714
$ws.Range("I107").Value = 570
$ws.Range("K107").Value = 570
$ws.Range("M107").Value = 1350
$ws.Range("H116").Value = 6233.25
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 6233.25
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 6233.25
$ws.Range("M116").Value = $null
$ws.Range("N116").Value = -13117.25
$ws.Range("H137").Value = 1551.3
$ws.Range("I137").Value = 1551.3
$ws.Range("K137").Value = 4653.9
$ws.Range("M137").Value = -2103.9
$ws.Range("H138").Value = 1825.9565
$ws.Range("I138").Value = 665.6667
$ws.Range("K138").Value = 1997.0001
$ws.Range("M138").Value = 3142.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6316
$ws.Range("I32").Value = 5286.2607
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 5286.2607
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -4999.2607
$ws.Range("N32").Value = -30574
$ws.Range("H132").Value = 4998.5
$ws.Range("I132").Value = 4997
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 14991
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -12461
$ws.Range("N132").Value = -20060

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = $null
$ws.Range("H86").Value = 1904.4333
$ws.Range("I86").Value = 1568.4117
$ws.Range("J86").Value = 2343.8462
$ws.Range("K86").Value = 1568.4117
$ws.Range("L86").Value = 2343.8462
$ws.Range("M86").Value = -445.4117000000001
$ws.Range("N86").Value = -4589.8462
$ws.Range("H89").Value = 1904.4333
$ws.Range("I89").Value = 1568.4117
$ws.Range("J89").Value = 2343.8462
$ws.Range("K89").Value = 7842.058500000001
$ws.Range("L89").Value = 11719.231
$ws.Range("M89").Value = -2226.058500000001
$ws.Range("N89").Value = -22951.231
$ws.Range("H94").Value = 2543.8
$ws.Range("I94").Value = 2543.8
$ws.Range("K94").Value = 2543.8
$ws.Range("M94").Value = -2092.8
$ws.Range("H140").Value = 148998
$ws.Range("J140").Value = 148998
$ws.Range("L140").Value = 148998
$ws.Range("N140").Value = -159358

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1587.25
$ws.Range("I107").Value = 1587.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1587.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 332.75
$ws.Range("N107").Value = $null
$ws.Range("H132").Value = 1100
$ws.Range("I132").Value = 1100
$ws.Range("K132").Value = 3300
$ws.Range("M132").Value = -770

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("H32").Value = 4950
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = $null
$ws.Range("H87").Value = 475
$ws.Range("I87").Value = 475
$ws.Range("K87").Value = 1425
$ws.Range("M87").Value = -177
$ws.Range("H90").Value = 475
$ws.Range("I90").Value = 475
$ws.Range("K90").Value = 4275
$ws.Range("M90").Value = 1965
$ws.Range("H129").Value = 359
$ws.Range("I129").Value = 100
$ws.Range("J129").Value = 488.5
$ws.Range("K129").Value = 300
$ws.Range("L129").Value = 1465.5
$ws.Range("M129").Value = 4700
$ws.Range("N129").Value = -11465.5
$ws.Range("H137").Value = 2729.75
$ws.Range("I137").Value = 2421.2
$ws.Range("J137").Value = 3244
$ws.Range("K137").Value = 7263.599999999999
$ws.Range("L137").Value = 9732
$ws.Range("M137").Value = -2163.599999999999
$ws.Range("N137").Value = -19932
$ws.Range("H140").Value = 2429.9092
$ws.Range("I140").Value = 2429.9092
$ws.Range("K140").Value = 7289.7276
$ws.Range("M140").Value = -2109.7276

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 36444.25
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 36444.25
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 36444.25
$ws.Range("M43").Value = $null
$ws.Range("N43").Value = -36746.25
$ws.Range("H47").Value = 35000
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -36136
$ws.Range("H55").Value = 25154.6
$ws.Range("I55").Value = 6500
$ws.Range("J55").Value = 37591
$ws.Range("K55").Value = 6500
$ws.Range("L55").Value = 37591
$ws.Range("M55").Value = -6173
$ws.Range("N55").Value = -38245
$ws.Range("H70").Value = 9999
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = $null
$ws.Range("H73").Value = 9999
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = $null
$ws.Range("H80").Value = 4796.3335
$ws.Range("I80").Value = 2399
$ws.Range("J80").Value = 5995
$ws.Range("K80").Value = 2399
$ws.Range("L80").Value = 5995
$ws.Range("M80").Value = -1401
$ws.Range("N80").Value = -7991
$ws.Range("H83").Value = 4796.3335
$ws.Range("I83").Value = 2399
$ws.Range("J83").Value = 5995
$ws.Range("K83").Value = 11995
$ws.Range("L83").Value = 29975
$ws.Range("M83").Value = -7003
$ws.Range("N83").Value = -39959
$ws.Range("H97").Value = 1500
$ws.Range("J97").Value = 1500
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492
$ws.Range("H122").Value = 8932190
$ws.Range("I122").Value = 10420222
$ws.Range("K122").Value = 31260666
$ws.Range("M122").Value = -31258216

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2380.4546
$ws.Range("I7").Value = 2154.5557
$ws.Range("K7").Value = 2154.5557
$ws.Range("M7").Value = -2042.5557
$ws.Range("H40").Value = 9793.714
$ws.Range("I40").Value = 9277.5
$ws.Range("J40").Value = 10000.2
$ws.Range("K40").Value = 9277.5
$ws.Range("L40").Value = 10000.2
$ws.Range("M40").Value = -9141.5
$ws.Range("N40").Value = -10272.2
$ws.Range("H82").Value = 1067.8572
$ws.Range("J82").Value = 1307.7
$ws.Range("L82").Value = 1307.7
$ws.Range("N82").Value = -2029.7
$ws.Range("H85").Value = 1067.8572
$ws.Range("J85").Value = 1307.7
$ws.Range("L85").Value = 1307.7
$ws.Range("N85").Value = -3803.7
$ws.Range("H126").Value = 2380.4546
$ws.Range("I126").Value = 2154.5557
$ws.Range("K126").Value = 6463.6671
$ws.Range("M126").Value = -3993.6671
$ws.Range("H132").Value = 10068.857
$ws.Range("I132").Value = 10913.667
$ws.Range("K132").Value = 32741.001
$ws.Range("M132").Value = -30211.001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7999.5
$ws.Range("J62").Value = 7999.5
$ws.Range("L62").Value = 7999.5
$ws.Range("N62").Value = -9247.5
$ws.Range("H65").Value = 7999.5
$ws.Range("J65").Value = 7999.5
$ws.Range("L65").Value = 39997.5
$ws.Range("N65").Value = -46237.5
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null
